# Refresh the crypto price/volume table: Tue Jul 23 05:19:15 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.832.80'
$ws.Range('E2').Value = '  -1.44%  '

# Row 3
$ws.Range('D3').Value = '3.451.63'
$ws.Range('E3').Value = '  -1.69%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.81%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.21%  '

# Row 7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.605'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.77%  '

# Row 9
$ws.Range('D9').Value = '3.450.98'
$ws.Range('E9').Value = '  -1.72%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.135'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.09%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.91'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.25%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.422'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.99%  '

# Row 13
$ws.Range('D13').Value = '4.048.48'
$ws.Range('E13').Value = '  -1.62%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.98'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.18%  '

# Row 15
$ws.Range('E15').Value = '  -3.29%  '

# Row 16
$ws.Range('D16').Value = '66.858.08'
$ws.Range('E16').Value = '  -1.38%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000173'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.90%  '

# Row 18
$ws.Range('D18').Value = '3.455.53'
$ws.Range('E18').Value = '  -1.57%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.56%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.83'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.32%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '377.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.22%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.81'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.77%  '

# Row 23
$ws.Range('E23').Value = '  +0.05%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.52%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '71.33'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.01%  '

# Row 26
$ws.Range('E26').Value = '  -2.14%  '

# Row 27
$ws.Range('E27').Value = '  -2.17%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.88'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.50%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.173'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.38%  '

# Row 30
$ws.Range('E30').Value = '  +0.35%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.91'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.42%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.01'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.68%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.88'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.16%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.35'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.33%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.15'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.50%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.53'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.49%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.33'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.33%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.878'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.09%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '27.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.97%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.98%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.62'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.17%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.60'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.10%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.12%  '

# Row 45
$ws.Range('D45').Value = '2.688.68'
$ws.Range('E45').Value = '  -6.41%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0694'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.25%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '25.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.77%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.96'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.05%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0293'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.11%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '319.92'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.23%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.32%  '
